$d = $word.ActiveDocument
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    Write-Output "$i|$t"
}
